# Refresh Leve profit-calculation market data across the Sheets workbook.
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

# ALC!row17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1384.6552
$ws.Range("J17").Value = 1498.2693
$ws.Range("L17").Value = 4494.8079
$ws.Range("N17").Value = -4830.8079

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5007182
$ws.Range("I138").Value = 1764353.6
$ws.Range("J138").Value = 6292076.5
$ws.Range("K138").Value = 5293060.800000001
$ws.Range("L138").Value = 18876229.5
$ws.Range("M138").Value = -5287920.800000001
$ws.Range("N138").Value = -18886509.5

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8639.613
$ws.Range("I74").Value = 1611.8125
$ws.Range("J74").Value = 16135.934
$ws.Range("K74").Value = 1611.8125
$ws.Range("L74").Value = 16135.934
$ws.Range("M74").Value = -737.8125
$ws.Range("N74").Value = -17883.934

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8639.613
$ws.Range("I77").Value = 1611.8125
$ws.Range("J77").Value = 16135.934
$ws.Range("K77").Value = 8059.0625
$ws.Range("L77").Value = 80679.67
$ws.Range("M77").Value = -3691.0625
$ws.Range("N77").Value = -89415.67

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2083.2173
$ws.Range("I122").Value = 1791.8334
$ws.Range("J122").Value = 2401.0908
$ws.Range("K122").Value = 5375.5002
$ws.Range("L122").Value = 7203.2724
$ws.Range("M122").Value = -2925.5002
$ws.Range("N122").Value = -12103.2724

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1115.9333
$ws.Range("I94").Value = 1354.8889
$ws.Range("J94").Value = 757.5
$ws.Range("K94").Value = 1354.8889
$ws.Range("L94").Value = 757.5
$ws.Range("M94").Value = -903.8888999999999
$ws.Range("N94").Value = -1659.5

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1426.9375
$ws.Range("I99").Value = 1161.9166
$ws.Range("K99").Value = 1161.9166
$ws.Range("M99").Value = 336.0834

# BSM!row115
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 36000
$ws.Range("J115").Value = 36000
$ws.Range("L115").Value = 36000
$ws.Range("N115").Value = -39134

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6945999.5
$ws.Range("I99").Value = 12501200
$ws.Range("J99").Value = 1998.75
$ws.Range("K99").Value = 12501200
$ws.Range("L99").Value = 1998.75
$ws.Range("M99").Value = -12499702
$ws.Range("N99").Value = -4994.75

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1771.1765
$ws.Range("I122").Value = 967.4167
$ws.Range("J122").Value = 3700.2
$ws.Range("K122").Value = 2902.2501
$ws.Range("L122").Value = 11100.6
$ws.Range("M122").Value = -452.2501000000002
$ws.Range("N122").Value = -16000.6

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6945999.5
$ws.Range("I126").Value = 12501200
$ws.Range("J126").Value = 1998.75
$ws.Range("K126").Value = 37503600
$ws.Range("L126").Value = 5996.25
$ws.Range("M126").Value = -37501130
$ws.Range("N126").Value = -10936.25

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2919.8965
$ws.Range("I134").Value = 1502.05
$ws.Range("J134").Value = 6070.6665
$ws.Range("K134").Value = 4506.15
$ws.Range("L134").Value = 18211.9995
$ws.Range("M134").Value = -1971.15
$ws.Range("N134").Value = -23281.9995

# CRP!row141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 264609.28
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 264609.28
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 264609.28
$ws.Range("N141").Value = -274969.28
$ws.Range("M141").ClearContents()

# CUL!row125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 3001.9412
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 3002.0625
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 9006.1875
$ws.Range("M125").Value = -4080
$ws.Range("N125").Value = -18846.1875

# CUL!row129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1504.5264
$ws.Range("I129").Value = 412.2857
$ws.Range("J129").Value = 2141.6667
$ws.Range("K129").Value = 1236.8571
$ws.Range("L129").Value = 6425.000100000001
$ws.Range("M129").Value = 3763.1429
$ws.Range("N129").Value = -16425.0001

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5552.2896
$ws.Range("I70").Value = 5440.794
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 5440.794
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -5170.794
$ws.Range("N70").Value = -7040

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5552.2896
$ws.Range("I73").Value = 5440.794
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 5440.794
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -4504.794
$ws.Range("N73").Value = -8372

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1816.7646
$ws.Range("I102").Value = 1527.5
$ws.Range("J102").Value = 3166.6667
$ws.Range("K102").Value = 1527.5
$ws.Range("L102").Value = 3166.6667
$ws.Range("M102").Value = 94.5
$ws.Range("N102").Value = -6410.6667

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 795331
$ws.Range("I122").Value = 2223563
$ws.Range("J122").Value = 1868.7778
$ws.Range("K122").Value = 6670689
$ws.Range("L122").Value = 5606.3334
$ws.Range("M122").Value = -6668239
$ws.Range("N122").Value = -10506.3334

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4740.3
$ws.Range("I132").Value = 4924.8887
$ws.Range("J132").Value = 4589.273
$ws.Range("K132").Value = 14774.6661
$ws.Range("L132").Value = 13767.819
$ws.Range("M132").Value = -12244.6661
$ws.Range("N132").Value = -18827.819

# LTW!row109
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 27000
$ws.Range("J109").Value = 27000
$ws.Range("L109").Value = 27000
$ws.Range("N109").Value = -29774

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3347.6086
$ws.Range("I122").Value = 1856.4286
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5569.2858
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3119.2858
$ws.Range("N122").Value = -16900

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5777.2
$ws.Range("I132").Value = 5307.778
$ws.Range("J132").Value = 6161.273
$ws.Range("K132").Value = 15923.334
$ws.Range("L132").Value = 18483.819
$ws.Range("M132").Value = -13393.334
$ws.Range("N132").Value = -23543.819

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5371.5186
$ws.Range("I136").Value = 2735.9285
$ws.Range("J136").Value = 8209.846
$ws.Range("K136").Value = 8207.7855
$ws.Range("L136").Value = 24629.538
$ws.Range("M136").Value = -5657.7855
$ws.Range("N136").Value = -29729.538

# WVR!row26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 19250
$ws.Range("I26").Value = 28500
$ws.Range("K26").Value = 28500
$ws.Range("M26").Value = -28207

# WVR!row62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7705328.5
$ws.Range("I62").Value = 11126518
$ws.Range("J62").Value = 7652.75
$ws.Range("K62").Value = 11126518
$ws.Range("L62").Value = 7652.75
$ws.Range("M62").Value = -11125894
$ws.Range("N62").Value = -8900.75

# WVR!row65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7705328.5
$ws.Range("I65").Value = 11126518
$ws.Range("J65").Value = 7652.75
$ws.Range("K65").Value = 55632590
$ws.Range("L65").Value = 38263.75
$ws.Range("M65").Value = -55629470
$ws.Range("N65").Value = -44503.75

# WVR!row69
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# WVR!row72
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 23593.762
$ws.Range("I122").Value = 37212.535
$ws.Range("J122").Value = 2409
$ws.Range("K122").Value = 111637.605
$ws.Range("L122").Value = 7227
$ws.Range("M122").Value = -109187.605
$ws.Range("N122").Value = -12127

# WVR!row135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 131407.5
$ws.Range("I135").Value = 55000
$ws.Range("K135").Value = 55000
$ws.Range("M135").Value = -49930

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15199264
$ws.Range("I136").Value = 30395768
$ws.Range("J136").Value = 2760.2727
$ws.Range("K136").Value = 91187304
$ws.Range("L136").Value = 8280.8181
$ws.Range("M136").Value = -91184754
$ws.Range("N136").Value = -13380.8181

# WVR!row141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 68571.664
$ws.Range("J141").Value = 68571.664
$ws.Range("L141").Value = 68571.664
$ws.Range("N141").Value = -78931.664
